$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 387). The value 45186 (2023-09-17) must be updated to
# 45188 (2023-09-19) for all of these rows.
$ws.Range("C2:C387").Value = 45188
